$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add new value on Sheet2!B1 so the multi-sheet range formulas on Sheet1
# (which span A1:B2 across Sheet1:Sheet3) pick up a fourth value.
$ws2.Range("B1").Value = 44

# Update the active cell selections to match the recorded workbook state.
$ws2.Range("B2").Select()
$ws1.Select()
$ws1.Range("I5").Select()

$wb.Save()
